# Update metrics values in Sheet1 (B2:Q26) to reflect the new LM training run.
# All data rows (2-26) previously shared identical values per column; the new
# values are likewise identical per column across all rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values per column (B..Q), applied uniformly to rows 2 through 26.
# (Written in plain decimal notation since the PowerShell parser here does
# not accept exponent literals like 1e-07.)
$newValues = @{
    "B" = 0.9999989698596451
    "C" = 0.9990244665669379
    "D" = 0.999995480989238
    "E" = 0.9999992516518608
    "F" = 0.9999972265768934
    "G" = 0.0000009615908033740045
    "H" = 0.0009106176387947784
    "I" = 0.000004614011809662338
    "J" = 0.0000005116618731179146
    "K" = 0.000002562836841390126
    "L" = 0.00005130350677945762
    "M" = 0.0009806073645318011
    "N" = 1.000024723368518
    "O" = 0.001022353850177004
    "P" = 77.70935367418031
    "Q" = 108.1812492958853
}

foreach ($row in 2..26) {
    foreach ($col in $newValues.Keys) {
        $ws.Range("$col$row").Value = $newValues[$col]
    }
}
